# Append a 1x3 "Name / age / work" table right after the final paragraph
# ("this is the third paragraph!!!"), before the section break — matching
# a table freshly inserted via Word's Insert > Table (default look/grid).

$d = $word.ActiveDocument

# Collapsed range at the very end of the document body (after the last
# paragraph's text, ahead of the section properties).
$endPos = $d.Content.End
$rng = $d.Range($endPos, $endPos)

$tableXml = '<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:tblPr>' +
        '<w:tblW w:type="auto" w:w="0"/>' +
        '<w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/>' +
    '</w:tblPr>' +
    '<w:tblGrid>' +
        '<w:gridCol w:w="2880"/>' +
        '<w:gridCol w:w="2880"/>' +
        '<w:gridCol w:w="2880"/>' +
    '</w:tblGrid>' +
    '<w:tr>' +
        '<w:tc>' +
            '<w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr>' +
            '<w:p><w:r><w:t>Name</w:t></w:r></w:p>' +
        '</w:tc>' +
        '<w:tc>' +
            '<w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr>' +
            '<w:p><w:r><w:t>age</w:t></w:r></w:p>' +
        '</w:tc>' +
        '<w:tc>' +
            '<w:tcPr><w:tcW w:type="dxa" w:w="2880"/></w:tcPr>' +
            '<w:p><w:r><w:t>work</w:t></w:r></w:p>' +
        '</w:tc>' +
    '</w:tr>' +
'</w:tbl>'

[void]$rng.InsertXML($tableXml)
